# Adds a new "param:scope" column to the Tests sheet, adds a new test row
# that exercises the missing-required-parameter case, and documents the new
# path parameter on the Documentation sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Tests"
# ---------------------------------------------------------------------
$tests = $wb.Worksheets.Item("Tests")

# New column I: header + width (~20 characters, matches the other columns'
# customWidth sizing)
$tests.Columns.Item(9).ColumnWidth = 19.17

function Set-Text($cell, $text) {
    # A leading apostrophe forces plain-text storage (matches the original
    # file's habit of storing every textual value as explicit text,
    # including ones that look like booleans/numbers, e.g. "true"/"400").
    # Resetting the style afterwards drops the resulting quote-prefix
    # formatting flag so the cell ends up as plain, unstyled text.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-Text $tests.Cells.Item(1, 9) "param:scope"

# Existing row 2 gets a value in the new column
Set-Text $tests.Cells.Item(2, 9) "work"

# New row 3 - "Missing Required Param" test
Set-Text $tests.Cells.Item(3, 1) "get-translate - Missing Required Param"
Set-Text $tests.Cells.Item(3, 2) "Test GET /api/translate/:scope with missing required parameters"
Set-Text $tests.Cells.Item(3, 3) "true"
$tests.Cells.Item(3, 4).Value = 400
$tests.Cells.Item(3, 5).Value = 10000
$tests.Cells.Item(3, 6).Value = 2000
$tests.Cells.Item(3, 7).Value = 500
Set-Text $tests.Cells.Item(3, 8) "get-translate,validation"
Set-Text $tests.Cells.Item(3, 9) ""

# ---------------------------------------------------------------------
# Sheet "Documentation"
# ---------------------------------------------------------------------
$docs = $wb.Worksheets.Item("Documentation")

# Capture the "Endpoint-Specific Notes" block text before overwriting
# anything, then shift it down by 3 rows to make room for the new
# "Parameter Descriptions" block.
$notesHeader = $docs.Cells.Item(18, 1).Value()
$methodLine = $docs.Cells.Item(19, 1).Value()
$pathLine = $docs.Cells.Item(20, 1).Value()
$descLine = $docs.Cells.Item(21, 1).Value()

Set-Text $docs.Cells.Item(24, 1) "• Required parameters: scope"
Set-Text $docs.Cells.Item(23, 1) $descLine
Set-Text $docs.Cells.Item(22, 1) $pathLine
Set-Text $docs.Cells.Item(21, 1) $methodLine
Set-Text $docs.Cells.Item(20, 1) $notesHeader
Set-Text $docs.Cells.Item(19, 1) ""

# New "Parameter Descriptions" block
Set-Text $docs.Cells.Item(17, 1) "Parameter Descriptions:"
Set-Text $docs.Cells.Item(18, 1) "param:scope"
Set-Text $docs.Cells.Item(18, 2) "Search scope (work, person, place, concept, event, etc.) (string) (REQUIRED - highlighted in yellow)"
